# edit.ps1 - applies the "almost finished part 1" revision to insurance_report.docx
#
# Section 1 (המבוטח -> כללי) and the start of Section 2 (מסוג האירוע -> התביעה)
# are rewritten from a field/value layout into free-text narrative
# paragraphs; one brand-new paragraph is added, and one paragraph
# (vehicle colour) is dropped entirely.
#
# Note: we deliberately locate text with Find.Execute and then assign the
# replacement directly to Range.Text (instead of passing it through
# Find.Execute's ReplaceWith argument) because the ReplaceWith path runs
# the typed "smart quotes" AutoCorrect pass and mangles the straight
# apostrophes ( ' ) that appear in the new Hebrew sentences into curly
# Unicode quotes.

$d = $word.ActiveDocument

function Set-ParaText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Set-ParaText: text not found: $oldText"
    }
    $rng.Text = $newText
}

# 1. "1. המבוטח" -> "1. כללי"
Set-ParaText '1. המבוטח' '1. כללי'

# 2. family-name line -> investigation-request narrative
Set-ParaText 'שם משפחה: אלכסנדר שניידרמן' 'נתבקשנו על ידי חברתכם לבצע חקירה בעקבות הודעת המבוטח על תאונה שארעה לו ברכבו מסוג: קיה בצבע פיקנטו כחול,משנת יצור 2020.'

# 3. ID-number line -> accident/vehicles-involved narrative
Set-ParaText 'מספר תעודת זהות: ' 'כתוצאה מהתאונה נפגעו שני הרכבים המעורבים – רכבו של המבוטח אופנוע צד ג'' מס'' רישוי:  554-49-103 מסוג ימהה טימקס בצבע כחול, 2019 משנת ייצור'

# 4. brand-new paragraph inserted right after the paragraph edited in step 3
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like '*כתוצאה מהתאונה נפגעו*') {
        $null = $p.Range.InsertParagraphAfter()
        $newP = $d.Paragraphs.Item($i + 1)
        $newP.Range.Text = 'להלן יובאו ממצאינו: '
        break
    }
}

# 5. "2. מסוג האירוע" -> "2. התביעה"
Set-ParaText '2. מסוג האירוע' '2. התביעה'

# 6. licence-number line -> claim narrative
Set-ParaText 'מספר רישוי רכב: 123-45-678' 'המדובר בתאונה בין המבוטח לאופנוע צד ג'''

# 7. event-date line -> insured-vehicle description
Set-ParaText 'תאריך האירוע: 31.12.2023' 'רכב המבוטח: מס'' רישוי 123-45-678מסוגקיהפיקנטובצבעכחולמנועבנזין,נפח מנוע 2000 סמ''''ק,  גיר אוטומטית, 4 דלתות, 7 כריות אוויר, שנת יצור 2020 .'

# 8. claim-number line -> vehicle ownership narrative
Set-ParaText 'מספר תביעה: 2418022441' 'הרכב רשום ע''''ש אלכסנדר שניידרמן, מר''חשנקין 13 ראשלצהרכב בבעלות רביעית'

# 9. vehicle-type line -> accident-date / test / license narrative
Set-ParaText 'סוג רכב: סדאן' 'התאונה דווחה שהתרחשה בתאריך ה- 31.12.2023תאריך טסט אחרון: 2023-01-01, תוקף רישיון2024-01-01 להלן צילום רישיון הרכב של המבוטח : '

# 10. drop the "vehicle colour" paragraph entirely
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like '*צבע רכב: כחול*') {
        $p.Range.Delete()
        break
    }
}

# 11. owner name/address line -> shortened address only
Set-ParaText 'שם וכתובת בעל הרכב: אלכסנדר שניידרמן, תל אביב' 'שם וכתובת בעל הרכב: שנקין 13 ראשלצ'
